$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: scroll position & selection ---
$ws.Range("I17").Select()

# --- Value changes ---
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 3
$ws.Range("F5").Value = 3
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 3
$ws.Range("F14").Value = 3
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 3
$ws.Range("F18").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 3
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 3
$ws.Range("F27").Value = 3
$ws.Range("F30").Value = 3
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 3
$ws.Range("F37").Value = 3
$ws.Range("F40").Value = 3
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = 3
$ws.Range("F42").Value = 3
$ws.Range("F45").Value = 3
$ws.Range("F51").Value = 3
$ws.Range("E52").Value = 0
$ws.Range("F52").Value = 3
$ws.Range("F53").Value = 3
$ws.Range("F54").Value = 3
$ws.Range("F55").Value = 3
$ws.Range("F56").Value = 3
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 3
$ws.Range("F65").Value = 3
$ws.Range("F67").Value = 3
$ws.Range("F69").Value = 3
$ws.Range("F71").Value = 3
$ws.Range("E72").Value = 0
$ws.Range("F72").Value = 3
$ws.Range("F75").Value = 3
$ws.Range("E80").Value = 0
$ws.Range("F80").Value = 3
$ws.Range("F81").Value = 3
$ws.Range("F87").Value = 3
$ws.Range("F88").Value = 3
$ws.Range("F90").Value = 3
$ws.Range("F91").Value = 3
$ws.Range("F92").Value = 3
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 3
$ws.Range("F94").Value = 3
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 3
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 3
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 3
$ws.Range("E100").Value = 0
$ws.Range("F100").Value = 3
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 3
$ws.Range("E102").Value = 0
$ws.Range("F102").Value = 3
$ws.Range("E103").Value = 0
$ws.Range("F103").Value = 3
$ws.Range("E104").Value = 0
$ws.Range("F104").Value = 3
$ws.Range("E105").Value = 0
$ws.Range("F105").Value = 3
$ws.Range("E106").Value = 0
$ws.Range("F106").Value = 3
$ws.Range("E107").Value = 0
$ws.Range("F107").Value = 3
$ws.Range("E108").Value = 0
$ws.Range("F108").Value = 3
$ws.Range("E109").Value = 0
$ws.Range("F109").Value = 3
$ws.Range("E111").Value = 0
$ws.Range("F111").Value = 3
$ws.Range("E112").Value = 0
$ws.Range("F112").Value = 3
$ws.Range("E113").Value = 0
$ws.Range("F113").Value = 3
$ws.Range("E117").Value = 0
$ws.Range("F117").Value = 3
$ws.Range("E118").Value = 0
$ws.Range("F118").Value = 3
$ws.Range("E121").Value = 0
$ws.Range("F121").Value = 3
$ws.Range("E122").Value = 0
$ws.Range("F122").Value = 3
$ws.Range("E123").Value = 0
$ws.Range("F123").Value = 3
$ws.Range("E124").Value = 0
$ws.Range("F124").Value = 3
$ws.Range("E128").Value = 0
$ws.Range("F128").Value = 3
$ws.Range("E129").Value = 0
$ws.Range("F129").Value = 3
$ws.Range("E130").Value = 0
$ws.Range("F130").Value = 3
$ws.Range("E131").Value = 0
$ws.Range("F131").Value = 3
$ws.Range("E132").Value = 0
$ws.Range("F132").Value = 3
$ws.Range("E135").Value = 0
$ws.Range("F135").Value = 3
$ws.Range("E136").Value = 0
$ws.Range("F136").Value = 3
$ws.Range("E137").Value = 0
$ws.Range("F137").Value = 3
$ws.Range("E138").Value = 0
$ws.Range("F138").Value = 3
$ws.Range("E139").Value = 0
$ws.Range("F139").Value = 3
$ws.Range("E142").Value = 0
$ws.Range("F142").Value = 3
$ws.Range("E146").Value = 0
$ws.Range("F146").Value = 3
$ws.Range("E147").Value = 0
$ws.Range("F147").Value = 3
$ws.Range("E148").Value = 0
$ws.Range("F148").Value = 3
$ws.Range("E149").Value = 0
$ws.Range("F149").Value = 3
$ws.Range("E150").Value = 0
$ws.Range("F150").Value = 3
$ws.Range("E154").Value = 0
$ws.Range("F154").Value = 3
$ws.Range("E155").Value = 0
$ws.Range("F155").Value = 3
$ws.Range("E157").Value = 0
$ws.Range("F157").Value = 3
$ws.Range("E159").Value = 0
$ws.Range("F159").Value = 3
$ws.Range("E160").Value = 0
$ws.Range("F160").Value = 3
$ws.Range("E161").Value = 0
$ws.Range("F161").Value = 3
$ws.Range("E162").Value = 0
$ws.Range("F162").Value = 3
$ws.Range("E166").Value = 0
$ws.Range("F166").Value = 3
$ws.Range("E168").Value = 0
$ws.Range("F168").Value = 3
$ws.Range("E177").Value = 0
$ws.Range("F177").Value = 3
$ws.Range("E183").Value = 0
$ws.Range("F183").Value = 3
$ws.Range("E184").Value = 0
$ws.Range("F184").Value = 3
$ws.Range("E187").Value = 0
$ws.Range("F187").Value = 3
$ws.Range("E188").Value = 0
$ws.Range("F188").Value = 3
$ws.Range("E191").Value = 0
$ws.Range("F191").Value = 3
$ws.Range("E192").Value = 0
$ws.Range("F192").Value = 3

# --- Style changes ---
$ws.Range("A13").Style = "Good"
$ws.Range("E13").Style = "Normal"
$ws.Range("A15").Style = "Good"
$ws.Range("E15").Style = "Normal"
$ws.Range("A20").Style = "Good"
$ws.Range("E20").Style = "Normal"
$ws.Range("A26").Style = "Good"
$ws.Range("E26").Style = "Normal"
$ws.Range("E35").Style = "Normal"
$ws.Range("A36").Style = "Good"
$ws.Range("E36").Style = "Normal"
$ws.Range("E43").Style = "Normal"
$ws.Range("A52").Style = "Good"
$ws.Range("A58").Style = "Good"
$ws.Range("A72").Style = "Good"
$ws.Range("A75").Style = "Good"
$ws.Range("E78").Style = "Normal"
$ws.Range("A80").Style = "Good"
$ws.Range("A96").Style = "Good"
$ws.Range("E96").Style = "Normal"
$ws.Range("A118").Style = "Good"
$ws.Range("E119").Style = "Normal"
$ws.Range("A137").Style = "Good"
$ws.Range("A146").Style = "Good"
$ws.Range("A157").Style = "Good"
$ws.Range("A183").Style = "Neutral"
